# Auto-generated PowerShell COM-interop script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update date (column B) for rows 1014-1019: 46005 -> 46007
for ($r = 1014; $r -le 1019; $r++) {
    $ws.Cells.Item($r, 2).Value = 46007
}

# 2) Append 13 new rows (1020-1032) by copying the formatting/style of row 1019
$srcRow = $ws.Range("A1019:V1019")
for ($r = 1020; $r -le 1032; $r++) {
    $dstRow = $ws.Range("A" + $r + ":V" + $r)
    $srcRow.Copy($dstRow)
}

# 3) Fill in the new rows' data
$rows = @(
    @{ Row=1020; E='Sofiane Belle'; F='left forward'; G='01:16:02'; H=5.56; I=0.69; J=4.86; K=0.47; L=0.19; M=0.04; N=0; O=4; P=4.25; Q=27.93; R=4.12; S=19; T=1; U=17; V=6 },
    @{ Row=1021; E='Emmanuel Valey'; F='left forward'; G='01:16:21'; H=6.34; I=1.16; J=5.16; K=0.65; L=0.37; M=0.15; N=0; O=14; P=4.46; Q=29.23; R=4.78; S=38; T=19; U=30; V=20 },
    @{ Row=1022; E='Ilyes Boughanmi'; F='center forward'; G='01:15:29'; H=4.59; I=0.33; J=4.26; K=0.24; L=0.09; M=0.01; N=0; O=1; P=3.19; Q=26.53; R=4.38; S=20; T=7; U=22; V=6 },
    @{ Row=1023; E='Ilan Ihaddadene'; F='center midfield'; G='01:15:16'; H=5.57; I=0.67; J=4.89; K=0.51; L=0.14; M=0.02; N=0; O=1; P=4.38; Q=26.27; R=4.47; S=25; T=3; U=5; V=3 },
    @{ Row=1024; E='Levy Ndoutoume'; F='left back'; G='01:15:17'; H=3.94; I=0.39; J=3.54; K=0.3; L=0.07; M=0.03; N=0; O=2; P=3.05; Q=28.34; R=4.66; S=11; T=5; U=6; V=2 },
    @{ Row=1025; E='Kamal Bafounta'; F='center midfield'; G='01:14:31'; H=5.51; I=0.62; J=4.88; K=0.51; L=0.09; M=0.03; N=0; O=3; P=4.36; Q=26.7; R=4.13; S=25; T=3; U=17; V=2 },
    @{ Row=1026; E='Yoan Zouma'; F='center back'; G='01:14:37'; H=4.54; I=0.35; J=4.18; K=0.29; L=0.07; M=0; N=0; O=0; P=3.56; Q=23.61; R=3.94; S=16; T=0; U=8; V=1 },
    @{ Row=1027; E='Karahali Souaré'; F='right forward'; G='01:13:58'; H=4.82; I=0.48; J=4.33; K=0.37; L=0.09; M=0.03; N=0; O=3; P=3.69; Q=29.74; R=5.01; S=31; T=13; U=29; V=4 },
    @{ Row=1028; E='Omar Benyounes'; F='center midfield'; G='01:16:15'; H=5.51; I=0.73; J=4.77; K=0.37; L=0.27; M=0.1; N=0; O=10; P=3.94; Q=29.86; R=5.39; S=42; T=18; U=36; V=12 },
    @{ Row=1029; E='Jeremie Laurent'; F='left forward'; G='01:14:49'; H=6.18; I=1.33; J=4.83; K=0.75; L=0.44; M=0.15; N=0; O=10; P=4.89; Q=30.13; R=5.28; S=39; T=15; U=33; V=19 },
    @{ Row=1030; E='Romain Thunet'; F='center back'; G='01:13:18'; H=4.88; I=0.34; J=4.53; K=0.24; L=0.09; M=0.02; N=0; O=1; P=3.95; Q=26.51; R=4.29; S=16; T=4; U=8; V=1 },
    @{ Row=1031; E='Malik Boussaid'; F='right back'; G='01:14:37'; H=4.76; I=0.5; J=4.25; K=0.29; L=0.18; M=0.04; N=0; O=4; P=3.4; Q=27.28; R=3.71; S=30; T=0; U=18; V=7 },
    @{ Row=1032; E='Mattheo Haon'; F='right back'; G='01:16:48'; H=6.0; I=0.77; J=5.21; K=0.52; L=0.21; M=0.05; N=0; O=3; P=4.61; Q=28.4; R=4.12; S=20; T=2; U=17; V=2 }
)

foreach ($d in $rows) {
    $r = $d.Row
    $ws.Cells.Item($r, 1).Value = "Entrainement"
    $ws.Cells.Item($r, 2).Value = 46008
    $ws.Cells.Item($r, 3).Value = "Global"
    $ws.Cells.Item($r, 4).Value = "J-3"
    $ws.Cells.Item($r, 5).Value = $d.E
    $ws.Cells.Item($r, 6).Value = $d.F
    $ws.Cells.Item($r, 7).Value = $d.G
    $ws.Cells.Item($r, 8).Value = $d.H
    $ws.Cells.Item($r, 9).Value = $d.I
    $ws.Cells.Item($r, 10).Value = $d.J
    $ws.Cells.Item($r, 11).Value = $d.K
    $ws.Cells.Item($r, 12).Value = $d.L
    $ws.Cells.Item($r, 13).Value = $d.M
    $ws.Cells.Item($r, 14).Value = $d.N
    $ws.Cells.Item($r, 15).Value = $d.O
    $ws.Cells.Item($r, 16).Value = $d.P
    $ws.Cells.Item($r, 17).Value = $d.Q
    $ws.Cells.Item($r, 18).Value = $d.R
    $ws.Cells.Item($r, 19).Value = $d.S
    $ws.Cells.Item($r, 20).Value = $d.T
    $ws.Cells.Item($r, 21).Value = $d.U
    $ws.Cells.Item($r, 22).Value = $d.V
}

# 4) Update the view selection to match the new extent
$ws.Range("C1036").Select()

